# Regenerate orders with updated distance/sizes.
# Simple global token substitution across the used range's text values:
#   D64 -> D69, D80 -> D86, D51 -> D55, S30 -> S31
# (S20/S25 and every other token are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Order matters only in that each pair is independent (no overlap in
# source/target tokens), so a straightforward sequential Replace is safe.
# LookAt:=2 (xlPart) so the token is matched as a substring (e.g. inside
# "Face13_D64_S25" or "Face13_D64_S25_l.png"), not only whole-cell matches.
$xlPart = 2
$used.Replace("D64", "D69", $xlPart) | Out-Null
$used.Replace("D80", "D86", $xlPart) | Out-Null
$used.Replace("D51", "D55", $xlPart) | Out-Null
$used.Replace("S30", "S31", $xlPart) | Out-Null
